# Hides the API key: replaces the "Diversificado_recs_random AWM" result
# table (rows 31-40 + NDCG label in row 41) with a different randomized
# recommendation run (different POIs / ids / scores), per the commit
# "escondendo a chave API".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 31 (index 0)
$ws.Range("B31").Value = 56
$ws.Range("C31").Value = "Sal Marinho"
$ws.Range("D31").Value = "bar"
$ws.Range("E31").Value = 1
$ws.Range("F31").Value = 0.526
$ws.Range("G31").Value = -12.999919
$ws.Range("H31").Value = -38.456829

# Row 32 (index 1)
$ws.Range("B32").Value = 52
$ws.Range("C32").Value = "Bar e Restaurante"
$ws.Range("D32").Value = "bar"
$ws.Range("E32").Value = 1
$ws.Range("F32").Value = 0.517
$ws.Range("G32").Value = -12.9125452
$ws.Range("H32").Value = -38.4971898

# Row 33 (index 2)
$ws.Range("B33").Value = 367
$ws.Range("C33").Value = "Ponta do Humaita"
$ws.Range("D33").Value = "tourist_attraction"
$ws.Range("E33").Value = 1
$ws.Range("F33").Value = 0.507
$ws.Range("G33").Value = -12.9299294
$ws.Range("H33").Value = -38.5351303

# Row 34 (index 3)
$ws.Range("B34").Value = 111
$ws.Range("C34").Value = "Praia Pedra do Sal"
$ws.Range("D34").Value = "beach"
$ws.Range("E34").Value = 1.146865557795592
$ws.Range("F34").Value = 0.573
$ws.Range("G34").Value = -12.9524313
$ws.Range("H34").Value = -38.3460781

# Row 35 (index 4)
$ws.Range("B35").Value = 72
$ws.Range("C35").Value = "Sentollas Bar e Restaurante"
$ws.Range("D35").Value = "bar"
$ws.Range("E35").Value = 1.127132897628926
$ws.Range("F35").Value = 0.5639999999999999
$ws.Range("G35").Value = -12.967751
$ws.Range("H35").Value = -38.4082736

# Row 36 (index 5)
$ws.Range("B36").Value = 47
$ws.Range("C36").Value = "Rua 15 Restaurante e Bar"
$ws.Range("D36").Value = "bar"
$ws.Range("E36").Value = 1
$ws.Range("F36").Value = 0.517
$ws.Range("G36").Value = -12.9795299
$ws.Range("H36").Value = -38.4299455

# Row 37 (index 6)
$ws.Range("B37").Value = 9
$ws.Range("C37").Value = "Galeria Canizares UFBA"
$ws.Range("D37").Value = "art_gallery"
$ws.Range("E37").Value = 0
$ws.Range("F37").Value = 0.016
$ws.Range("G37").Value = -12.9911905
$ws.Range("H37").Value = -38.5211528

# Row 38 (index 7)
$ws.Range("B38").Value = 105
$ws.Range("C38").Value = "Praia de Placafor"
$ws.Range("D38").Value = "beach"
$ws.Range("E38").Value = 1.219282701737286
$ws.Range("F38").Value = 0.61
$ws.Range("G38").Value = -12.9512414
$ws.Range("H38").Value = -38.371031

# Row 39 (index 8)
$ws.Range("B39").Value = 82
$ws.Range("C39").Value = "Praia de Piata"
$ws.Range("D39").Value = "beach"
$ws.Range("E39").Value = 1
$ws.Range("F39").Value = 0.538
$ws.Range("G39").Value = -12.9547946
$ws.Range("H39").Value = -38.3826836

# Row 40 (index 9)
$ws.Range("B40").Value = 61
$ws.Range("C40").Value = "Beckels Pizza"
$ws.Range("D40").Value = "bar"
$ws.Range("E40").Value = 1
$ws.Range("F40").Value = 0.517
$ws.Range("G40").Value = -12.9919109
$ws.Range("H40").Value = -38.4551933

# Row 41: updated NDCG score label for this table
$ws.Range("A41").Value = "NDCG: 0.9570643210427647"
